$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '60.739.83'
$ws.Range('E2').Value = '  -1.80%  '

# Row 3
$ws.Range('D3').Value = '3.384.41'
$ws.Range('E3').Value = '  -2.23%  '

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.03%  '

# Row 5
$ws.Range('E5').Value = '  -2.32%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '141.88'
$ws.Range('E6').Value = '  -3.80%  '

# Row 7
$ws.Range('E7').Value = '  +0.06%  '

# Row 8
$ws.Range('D8').Value = '3.384.74'
$ws.Range('E8').Value = '  -2.24%  '

# Row 9
$ws.Range('E9').Value = '  -0.30%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '7.50'
$ws.Range('E10').Value = '  -2.23%  '

# Row 11
$ws.Range('E11').Value = '  -2.16%  '

# Row 12
$ws.Range('E12').Value = '  +2.31%  '

# Row 13
$ws.Range('D13').Value = '3.961.94'
$ws.Range('E13').Value = '  -2.24%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '28.38'
$ws.Range('E14').Value = '  +1.41%  '

# Row 15
$ws.Range('E15').Value = '  +1.54%  '

# Row 16
$ws.Range('E16').Value = '  -2.23%  '

# Row 17
$ws.Range('D17').Value = '3.383.31'
$ws.Range('E17').Value = '  -2.42%  '

# Row 18
$ws.Range('D18').Value = '60.807.66'
$ws.Range('E18').Value = '  -1.83%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.26'

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.07'
$ws.Range('E20').Value = '  -2.21%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '9.03'
$ws.Range('E21').Value = '  -5.78%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '384.75'
$ws.Range('E22').Value = '  -1.24%  '

# Row 23
$ws.Range('E23').Value = '  -0.79%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '73.68'
$ws.Range('E24').Value = '  -0.07%  '

# Row 25
$ws.Range('E25').Value = '  +0.24%  '

# Row 26
$ws.Range('E26').Value = '  -5.67%  '

# Row 27
$ws.Range('D27').Value = '3.522.10'
$ws.Range('E27').Value = '  -2.19%  '

# Row 28
$ws.Range('E28').Value = '  -2.59%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.26%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.44'
$ws.Range('E30').Value = '  -3.78%  '

# Row 31
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '8.02'
$ws.Range('E31').Value = '  -2.58%  '

# Row 32
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.44'
$ws.Range('E32').Value = '  -2.63%  '

# Row 33
$ws.Range('E33').Value = '  -2.06%  '

# Row 34
$ws.Range('E34').Value = '  +0.00%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '23.70'
$ws.Range('E35').Value = '  -2.19%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '7.00'
$ws.Range('E36').Value = '  -0.27%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '166.64'
$ws.Range('E37').Value = '  -0.23%  '

# Row 38
$ws.Range('E38').Value = '  -3.06%  '

# Row 39
$ws.Range('D39').Value = '3.416.09'
$ws.Range('E39').Value = '  -2.04%  '

# Row 40
$ws.Range('E40').Value = '  -4.91%  '

# Row 41
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '27.90'
$ws.Range('E41').Value = '  +1.98%  '

# Row 42
$ws.Range('B42').Value = 'Hedera'
$ws.Range('C42').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0779'
$ws.Range('E42').Value = '  -0.72%  '

# Row 43
$ws.Range('E43').Value = '  -3.17%  '

# Row 44
$ws.Range('E44').Value = '  -0.10%  '

# Row 45
$ws.Range('E45').Value = '  -1.66%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '41.73'
$ws.Range('E46').Value = '  -1.97%  '

# Row 47
$ws.Range('E47').Value = '  -3.02%  '

# Row 48
$ws.Range('D48').Value = '2.534.62'
$ws.Range('E48').Value = '  -1.41%  '

# Row 49
$ws.Range('E49').Value = '  -4.39%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '23.58'
$ws.Range('E50').Value = '  +1.82%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '6.87'
$ws.Range('E51').Value = '  -1.18%  '
